$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full new data set (rows 2-31), in order, replacing/extending the old A2:C21 block.
$data = @(
    @(-7.313242021728997, -8.175233658622219, 2.331058922935893),
    @(-8.089916832306811, -7.794766384012551, 4.25856484151354),
    @(-8.562996188799541, -5.62869268538905, 4.021127045154572),
    @(-8.465083991779997, -5.572577813092403, 3.310186512330003),
    @(-7.515201767285651, -5.813795685768134, 2.080139199892673),
    @(-3.889640840829573, -6.183977889079673, 3.298035799288292),
    @(-2.032402010524974, -7.917088676901407, 5.917996553813698),
    @(-0.8607274410771413, -18.58178931591555, 10.12933701159906),
    @(-3.271800985523277, -15.36747401835878, 6.096988180104466),
    @(-14.65323523914104, -4.466919653555957, -0.7585973178638932),
    @(-22.51664300058401, -7.749673447188126, -5.398741133072749),
    @(-18.3769370060341, -24.21302471908869, 0.6640657630621334),
    @(-11.91456683944254, -23.90803286608531, 10.49638333040125),
    @(-4.779002170936772, -4.17006232224267, 11.643902759926),
    @(0.3725347098182305, -30.31529621750734, 8.55354468495236),
    @(-15.74266366397642, -20.48928019579699, 1.483265680425281),
    @(-19.33448225376665, 2.267985530928218, 5.825747468892169),
    @(-1.347560611425691, -14.70038807625862, 10.01127913886427),
    @(7.625993616440701, -13.20338362104757, 14.98549651398371),
    @(-0.7890216976987219, -0.09984105241075003, 12.88112398222394),
    @(-7.435094737539091, -33.73446517364671, 13.88054672058896),
    @(7.333738565444815, -26.17261409759575, 7.461308479309277),
    @(0.3924275192565041, 8.770900960062418, -10.52958610946054),
    @(-27.41900163538287, -21.30080885045613, 1.806023990406281),
    @(-26.291865657358, -17.62141748035735, 14.67751483356255),
    @(-12.81046040852863, 7.074593609454595, 12.13140960768156),
    @(-5.703531527051739, -11.14156565946689, 10.64924546316558),
    @(5.73159689061778, -12.49397951013913, 3.91989309647508),
    @(-10.9989599106356, 0.90863177355609, -9.010537198945535),
    @(9.205162721520956, -32.48499697329918, -12.15963486129167)
)

$row = 2
foreach ($values in $data) {
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $row++
}
